# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values, matching the GitHub Actions scrape
# commit "Updated cryptos list on Thu Aug 31 17:53:21 UTC 2023".
#
# Price values that look like plain numbers (e.g. "1.008", "210.90") are
# written with a leading apostrophe so Excel keeps them as text (matching
# the source data, which stores every Price/Volume cell as a string -
# some look like thousands-dotted numbers such as "26.523.40" and must
# not be reinterpreted as numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.523.40' },
    @{ Cell = 'E2'; Value = '  -2.82%  ' },
    @{ Cell = 'D3'; Value = '1.671.68' },
    @{ Cell = 'E3'; Value = '  -2.16%  ' },
    @{ Cell = 'E4'; Value = '  +0.40%  ' },
    @{ Cell = 'D5'; Value = '''219.53' },
    @{ Cell = 'E5'; Value = '  -1.81%  ' },
    @{ Cell = 'D6'; Value = '''0.5159' },
    @{ Cell = 'E6'; Value = '  -2.83%  ' },
    @{ Cell = 'D7'; Value = '''1.008' },
    @{ Cell = 'E7'; Value = '  +0.45%  ' },
    @{ Cell = 'D8'; Value = '''0.06475' },
    @{ Cell = 'E8'; Value = '  -1.80%  ' },
    @{ Cell = 'D9'; Value = '''0.2577' },
    @{ Cell = 'E9'; Value = '  -2.97%  ' },
    @{ Cell = 'E10'; Value = '  -4.32%  ' },
    @{ Cell = 'D11'; Value = '''0.07672' },
    @{ Cell = 'E11'; Value = '  +0.33%  ' },
    @{ Cell = 'D12'; Value = '1.681.94' },
    @{ Cell = 'E12'; Value = '  -1.67%  ' },
    @{ Cell = 'D13'; Value = '''4.349' },
    @{ Cell = 'E13'; Value = '  -5.08%  ' },
    @{ Cell = 'D14'; Value = '1.901.10' },
    @{ Cell = 'E14'; Value = '  -2.21%  ' },
    @{ Cell = 'D15'; Value = '''0.5575' },
    @{ Cell = 'E15'; Value = '  -2.85%  ' },
    @{ Cell = 'D16'; Value = '0.0₅8045' },
    @{ Cell = 'E16'; Value = '  -1.83%  ' },
    @{ Cell = 'D17'; Value = '''64.91' },
    @{ Cell = 'E17'; Value = '  -3.94%  ' },
    @{ Cell = 'D18'; Value = '26.562.08' },
    @{ Cell = 'E18'; Value = '  -2.69%  ' },
    @{ Cell = 'D19'; Value = '''1.008' },
    @{ Cell = 'E19'; Value = '  +0.43%  ' },
    @{ Cell = 'D20'; Value = '''210.90' },
    @{ Cell = 'E20'; Value = '  -2.55%  ' },
    @{ Cell = 'D21'; Value = '''4.448' },
    @{ Cell = 'E21'; Value = '  -4.91%  ' },
    @{ Cell = 'D22'; Value = '''10.12' },
    @{ Cell = 'E22'; Value = '  -2.81%  ' },
    @{ Cell = 'D23'; Value = '''5.902' },
    @{ Cell = 'E23'; Value = '  -1.22%  ' },
    @{ Cell = 'D24'; Value = '''1.008' },
    @{ Cell = 'E24'; Value = '  +0.41%  ' },
    @{ Cell = 'D25'; Value = '''143.18' },
    @{ Cell = 'E25'; Value = '  +0.77%  ' },
    @{ Cell = 'D26'; Value = '''1.715' },
    @{ Cell = 'E26'; Value = '  -1.80%  ' },
    @{ Cell = 'D27'; Value = '''0.1170' },
    @{ Cell = 'E27'; Value = '  -3.89%  ' },
    @{ Cell = 'D28'; Value = '''6.996' },
    @{ Cell = 'E28'; Value = '  -3.65%  ' },
    @{ Cell = 'D29'; Value = '''15.77' },
    @{ Cell = 'E29'; Value = '  -3.55%  ' },
    @{ Cell = 'D30'; Value = '''0.05222' },
    @{ Cell = 'E30'; Value = '  -3.10%  ' },
    @{ Cell = 'D31'; Value = '''1.266' },
    @{ Cell = 'E31'; Value = '  -2.04%  ' },
    @{ Cell = 'E32'; Value = '  -4.43%  ' },
    @{ Cell = 'D33'; Value = '''3.212' },
    @{ Cell = 'E33'; Value = '  -6.21%  ' },
    @{ Cell = 'D34'; Value = '''1.581' },
    @{ Cell = 'E34'; Value = '  -3.77%  ' },
    @{ Cell = 'D35'; Value = '''2.766' },
    @{ Cell = 'E35'; Value = '  -4.00%  ' },
    @{ Cell = 'D36'; Value = '''2.377' },
    @{ Cell = 'E36'; Value = '  -1.89%  ' },
    @{ Cell = 'D37'; Value = '''0.9272' },
    @{ Cell = 'E37'; Value = '  -2.20%  ' },
    @{ Cell = 'D38'; Value = '''0.5731' },
    @{ Cell = 'E38'; Value = '  -2.14%  ' },
    @{ Cell = 'D39'; Value = '1.153.77' },
    @{ Cell = 'E39'; Value = '  +10.57%  ' },
    @{ Cell = 'D40'; Value = '''0.01591' },
    @{ Cell = 'E40'; Value = '  -2.67%  ' },
    @{ Cell = 'D41'; Value = '''1.008' },
    @{ Cell = 'E41'; Value = '  +0.40%  ' },
    @{ Cell = 'D42'; Value = '''0.8365' },
    @{ Cell = 'E42'; Value = '  -0.34%  ' },
    @{ Cell = 'D43'; Value = '''5.643' },
    @{ Cell = 'D44'; Value = '''100.01' },
    @{ Cell = 'E44'; Value = '  -1.07%  ' },
    @{ Cell = 'D45'; Value = '1.809.86' },
    @{ Cell = 'E45'; Value = '  -2.24%  ' },
    @{ Cell = 'D46'; Value = '0.0₈112' },
    @{ Cell = 'E46'; Value = '  -2.87%  ' },
    @{ Cell = 'D47'; Value = '''0.4495' },
    @{ Cell = 'E47'; Value = '  -0.14%  ' },
    @{ Cell = 'D48'; Value = '''55.77' },
    @{ Cell = 'E48'; Value = '  -4.05%  ' },
    @{ Cell = 'E49'; Value = '  -0.24%  ' },
    @{ Cell = 'D50'; Value = '''7.895' },
    @{ Cell = 'E50'; Value = '  -2.49%  ' },
    @{ Cell = 'D51'; Value = '''0.05137' },
    @{ Cell = 'E51'; Value = '  -1.97%  ' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
